# Update countries & provincias Spain
#
# The source data table ("Pais" sheet) is sorted by total cases (column B)
# descending. This refresh:
#   - updates the "last updated" timestamp (A1)
#   - updates several countries' case figures
#   - re-sorts a few rows whose totals changed enough to cross a
#     neighbouring country (Armenia now outranks Ghana/Suiza/Uzbekistan,
#     Hungria now outranks Guinea, Islas Malvinas now outranks Montserrat)
#
# Because the re-sort only ever swaps adjacent rows, each row is rewritten
# in place (label + all 7 data columns) rather than literally moving rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "last updated" timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 10 de Septiembre de 2020 a las 09:43"

# Row 7: Rusia - updated case figures
$ws.Cells.Item(7, 2).Value = 1046370
$ws.Cells.Item(7, 3).Value = 5363
$ws.Cells.Item(7, 4).Value = 862373
$ws.Cells.Item(7, 5).Value = 165734
$ws.Cells.Item(7, 7).Value = 128
$ws.Cells.Item(7, 8).Value = 18263

# Row 53: Barein - updated case figures
$ws.Cells.Item(53, 5).Value = 5427
$ws.Cells.Item(53, 7).Value = 1
$ws.Cells.Item(53, 8).Value = 204

# Rows 60-63: Armenia moves up ahead of Ghana/Suiza/Uzbekistan
# Row 60: was Ghana -> now Armenia (updated figures)
$ws.Cells.Item(60, 1).Value = "Armenia"
$ws.Cells.Item(60, 2).Value = 45326
$ws.Cells.Item(60, 3).Value = 174
$ws.Cells.Item(60, 4).Value = 41233
$ws.Cells.Item(60, 5).Value = 3187
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 906

# Row 61: was Suiza -> now Ghana (figures unchanged, just shifted down)
$ws.Cells.Item(61, 1).Value = "Ghana"
$ws.Cells.Item(61, 2).Value = 45313
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = 44188
$ws.Cells.Item(61, 5).Value = 842
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 283

# Row 62: was Uzbekistan -> now Suiza (figures unchanged, just shifted down)
$ws.Cells.Item(62, 1).Value = "Suiza"
$ws.Cells.Item(62, 2).Value = 45306
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 38100
$ws.Cells.Item(62, 5).Value = 5187
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 2019

# Row 63: was Armenia -> now Uzbekistan (figures unchanged, just shifted down)
$ws.Cells.Item(63, 1).Value = "Uzbekistan"
$ws.Cells.Item(63, 2).Value = 45160
$ws.Cells.Item(63, 3).Value = 230
$ws.Cells.Item(63, 4).Value = 42212
$ws.Cells.Item(63, 5).Value = 2580
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 2
$ws.Cells.Item(63, 8).Value = 368

# Row 66: Afganistan - updated case figures
$ws.Cells.Item(66, 2).Value = 38572
$ws.Cells.Item(66, 3).Value = 28
$ws.Cells.Item(66, 4).Value = 31129
$ws.Cells.Item(66, 5).Value = 6023

# Rows 95-96: Hungria moves up ahead of Guinea
# Row 95: was Guinea -> now Hungria (updated figures)
$ws.Cells.Item(95, 1).Value = "Hungria"
$ws.Cells.Item(95, 2).Value = 10191
$ws.Cells.Item(95, 3).Value = 476
$ws.Cells.Item(95, 4).Value = 3990
$ws.Cells.Item(95, 5).Value = 5571
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 2
$ws.Cells.Item(95, 8).Value = 630

# Row 96: was Hungria -> now Guinea (figures unchanged, just shifted down)
$ws.Cells.Item(96, 1).Value = "Guinea"
$ws.Cells.Item(96, 2).Value = 9885
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 9068
$ws.Cells.Item(96, 5).Value = 754
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 63

# Row 132: Lituania - updated case figures
$ws.Cells.Item(132, 2).Value = 3199
$ws.Cells.Item(132, 4).Value = 2030
$ws.Cells.Item(132, 5).Value = 1083

# Rows 214-215: Islas Malvinas moves up ahead of Montserrat
# Row 214: was Montserrat -> now Islas Malvinas (figures unchanged, just shifted up)
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 2).Value = 13
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 0

# Row 215: was Islas Malvinas -> now Montserrat (figures unchanged, just shifted down)
$ws.Cells.Item(215, 1).Value = "Montserrat"
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 1
